$wb = $excel.ActiveWorkbook

# --- Insert the new "Texas Notes" sheet between "Calculations" and "PEUDfSbQL" ---
$calcSheet = $wb.Worksheets.Item("Calculations")
$notesSheet = $wb.Worksheets.Add($null, $calcSheet)
$notesSheet.Name = "Texas Notes"

$notesSheet.Range("A1").Value  = "They are just comparing the efficiency of new appliances:"
$notesSheet.Range("A2").Value  = "standard versus energy star rebate qualifying"
$notesSheet.Range("A3").Value  = "i.e., the point of this spreadsheet is to estimate how much a household's energy consumption would change"
$notesSheet.Range("A4").Value  = "if they decide to use a rebate to get a higher-efficiency appliance instead of just opting for the cheap alternative. "
$notesSheet.Range("A6").Value  = "I think the method makes sense"
$notesSheet.Range("A8").Value  = "And there's no reason that Texas should be different. "
$notesSheet.Range("A9").Value  = "New technology in Texas should be as efficient as new technology across the US. "
$notesSheet.Range("A10").Value = "The only difference might be if Texas rebates incentivize a different level of efficiency than"
$notesSheet.Range("A11").Value = "national rebates do, but some of the other sources used in the building input files seem"
$notesSheet.Range("A12").Value = "to indicate that Texas doesn't usually have appliance rebates on top of the national ones. "
$notesSheet.Range("A14").Value = "So it's a good assumption that if a Texan uses a rebate to buy a more efficient appliance, that"
$notesSheet.Range("A15").Value = "rebate will be a national one and it will be based on national energy star standards."

# --- Update the stored selections on the other sheets ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("C6").Select()

$wsCalc = $wb.Worksheets.Item("Calculations")
$wsCalc.Activate()
$wsCalc.Range("B18").Select()

# --- Make "PEUDfSbQL" the active/selected sheet again, with its own selection ---
$wsPEU = $wb.Worksheets.Item("PEUDfSbQL")
$wsPEU.Activate()
$wsPEU.Range("B7").Select()
